$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.103.80"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -1.40%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.780.58"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -2.39%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'1.005"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  +0.14%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'336.39"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -2.77%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("E6").Value = "'  +0.13%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.3814"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -0.45%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  -3.70%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'48.08"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -2.61%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'1.184"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -4.50%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.07421"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -4.98%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'1.002"
$ws.Range("D12").Style = "Normal"
$ws.Range("E13").Value = "'  -3.16%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'6.435"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -3.36%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'1.781.65"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -2.04%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("E16").Value = "'  -3.34%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'0.00001082"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -4.35%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("E18").Value = "'  -1.69%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'83.23"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -4.09%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'  +0.16%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'6.528"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -0.75%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'  -2.95%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'27.122.99"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -1.38%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'12.19"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -8.67%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'2.374"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -3.40%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'2.496"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -7.80%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'1.453"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -4.29%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'21.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -6.31%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'155.18"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +0.79%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'1.980.38"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -2.03%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'133.73"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -2.35%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'3.973"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -0.84%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'6.009"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -6.74%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'0.08666"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -1.58%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'13.04"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -8.09%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'1.622"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -4.78%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'0.6816"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -4.03%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'5.369"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -5.78%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.06269"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -4.76%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  -5.59%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.02311"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -4.96%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'8.515"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -6.17%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'1.230"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -5.70%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'14.22"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -4.45%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("E46").Value = "'  -3.49%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'3.855"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -4.80%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'  -4.01%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'  -1.66%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.07077"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -4.03%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'78.51"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -3.19%  "
$ws.Range("E51").Style = "Normal"
